$wb = $excel.ActiveWorkbook

# Build the new "calc1"/"calc2" helper sheets (needed by the runorderall
# routine) by copying the existing "Sheet1"/"Sheet2" scratch sheets to the
# end of the tab order, then deleting the originals and renaming the
# copies. Doing it via Copy (rather than Move) mirrors how the sheetIds
# ended up allocated: calc1 -> 18, calc2 -> 19.

# Copy "Sheet1" to the very end of the workbook (Excel names the copy
# "Sheet1 (2)"). Worksheet handles in this host are position-bound, so each
# sheet is re-fetched by name right before it is used rather than re-using
# a reference captured earlier (stale references break once Delete()
# shifts the tab positions around).
$wb.Worksheets.Item("Sheet1").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Copy "Sheet2" so it lands immediately before the "Sheet1" copy.
$wb.Worksheets.Item("Sheet2").Copy($wb.Worksheets.Item("Sheet1 (2)"), $null)

# Remove the original scratch sheets, keeping only the copies at the end.
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

# Rename the copies to their final names.
$wb.Worksheets.Item("Sheet2 (2)").Name = "calc2"
$wb.Worksheets.Item("Sheet1 (2)").Name = "calc1"

# Re-enter the customer names on the expense_reports sheet so the report
# matches the freshly re-typed source values (clean strings, no stray
# trailing whitespace/newlines carried over from the old entries).
$expense = $wb.Worksheets.Item("expense_reports")
$expense.Range("B2").Value = "Carson Goble"
$expense.Range("B3").Value = "Aiden Herrera"
$expense.Range("B4").Value = "Cayden Doyle"

$expense.Activate()
